# Update TPM-derived NATMI metrics for Tnc-Itga7 LR pair (rows 2-10)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.1346003333333333
$ws.Range("H2").Value = 0.403801
$ws.Range("I2").Value = 0.009651054304565105
$ws.Range("J2").Value = 0.009651054304565105
$ws.Range("M2").Value = 1.743137
$ws.Range("N2").Value = 5.229411
$ws.Range("O2").Value = 0.03144673183548247
$ws.Range("P2").Value = 0.03144673183548247
$ws.Range("Q2").Value = 0.2346268212456667
$ws.Range("R2").Value = 2.111641391211
$ws.Range("S2").Value = 0.0003034941166453377
$ws.Range("T2").Value = 0.0003034941166453377

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.1346003333333333
$ws.Range("H3").Value = 0.403801
$ws.Range("I3").Value = 0.009651054304565105
$ws.Range("J3").Value = 0.009651054304565105
$ws.Range("N3").Value = 3.848628
$ws.Range("O3").Value = 0.02314348071905789
$ws.Range("P3").Value = 0.02314348071905789
$ws.Range("Q3").Value = 0.1726755372253334
$ws.Range("R3").Value = 1.554079835028
$ws.Range("S3").Value = 0.0002233589892162832
$ws.Range("T3").Value = 0.0002233589892162832

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.1346003333333333
$ws.Range("H4").Value = 0.403801
$ws.Range("I4").Value = 0.009651054304565105
$ws.Range("J4").Value = 0.009651054304565105
$ws.Range("M4").Value = 52.405407
$ws.Range("N4").Value = 157.216221
$ws.Range("O4").Value = 0.9454097874454597
$ws.Range("P4").Value = 0.9454097874454597
$ws.Range("Q4").Value = 7.053785250669002
$ws.Range("R4").Value = 63.48406725602101
$ws.Range("S4").Value = 0.009124201198703484
$ws.Range("T4").Value = 0.009124201198703484

# Row 5
$ws.Range("I5").Value = 0.8124788779145131
$ws.Range("J5").Value = 0.8124788779145132
$ws.Range("M5").Value = 1.743137
$ws.Range("N5").Value = 5.229411
$ws.Range("O5").Value = 0.03144673183548247
$ws.Range("P5").Value = 0.03144673183548247
$ws.Range("Q5").Value = 19.752177372389
$ws.Range("R5").Value = 177.769596351501
$ws.Range("S5").Value = 0.0255498053957714
$ws.Range("T5").Value = 0.0255498053957714

# Row 6
$ws.Range("I6").Value = 0.8124788779145131
$ws.Range("J6").Value = 0.8124788779145132
$ws.Range("N6").Value = 3.848628
$ws.Range("O6").Value = 0.02314348071905789
$ws.Range("P6").Value = 0.02314348071905789
$ws.Range("S6").Value = 0.01880358924565633
$ws.Range("T6").Value = 0.01880358924565633

# Row 7
$ws.Range("I7").Value = 0.8124788779145131
$ws.Range("J7").Value = 0.8124788779145132
$ws.Range("M7").Value = 52.405407
$ws.Range("N7").Value = 157.216221
$ws.Range("O7").Value = 0.9454097874454597
$ws.Range("P7").Value = 0.9454097874454597
$ws.Range("Q7").Value = 593.8264716635791
$ws.Range("R7").Value = 5344.438244972212
$ws.Range("S7").Value = 0.7681254832730854
$ws.Range("T7").Value = 0.7681254832730855

# Row 8
$ws.Range("G8").Value = 2.4807
$ws.Range("H8").Value = 7.4421
$ws.Range("I8").Value = 0.1778700677809217
$ws.Range("J8").Value = 0.1778700677809217
$ws.Range("M8").Value = 1.743137
$ws.Range("N8").Value = 5.229411
$ws.Range("O8").Value = 0.03144673183548247
$ws.Range("P8").Value = 0.03144673183548247
$ws.Range("Q8").Value = 4.3241999559
$ws.Range("R8").Value = 38.9177996031
$ws.Range("S8").Value = 0.005593432323065736
$ws.Range("T8").Value = 0.005593432323065736

# Row 9
$ws.Range("G9").Value = 2.4807
$ws.Range("H9").Value = 7.4421
$ws.Range("I9").Value = 0.1778700677809217
$ws.Range("J9").Value = 0.1778700677809217
$ws.Range("N9").Value = 3.848628
$ws.Range("O9").Value = 0.02314348071905789
$ws.Range("P9").Value = 0.02314348071905789
$ws.Range("Q9").Value = 3.1824304932
$ws.Range("R9").Value = 28.6418744388
$ws.Range("S9").Value = 0.004116532484185282
$ws.Range("T9").Value = 0.004116532484185281

# Row 10
$ws.Range("G10").Value = 2.4807
$ws.Range("H10").Value = 7.4421
$ws.Range("I10").Value = 0.1778700677809217
$ws.Range("J10").Value = 0.1778700677809217
$ws.Range("M10").Value = 52.405407
$ws.Range("N10").Value = 157.216221
$ws.Range("O10").Value = 0.9454097874454597
$ws.Range("P10").Value = 0.9454097874454597
$ws.Range("Q10").Value = 130.0020931449
$ws.Range("R10").Value = 1170.0188383041
$ws.Range("S10").Value = 0.1681601029736707
$ws.Range("T10").Value = 0.1681601029736707
